# PO Status update via batch
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last updated" timestamp banner in A1
$ws.Range("A1").Value = "Last updated: 2025-07-17 10:54:49"

# Row 7 - 4516351202_AIZU
$ws.Range("B7").Value = 19
$ws.Range("C7").Value = 0
$ws.Range("E7").Value = 4
$ws.Range("G7").Value = 14
$ws.Range("I7").Value = 0

# Row 8 - 4516351202_ARD
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 2
$ws.Range("G8").Value = 19
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 0

# Row 17 - 4516351202_TICL-FT
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").ClearContents()

# Row 18 - 4516351202_TICL-PR
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("G18").Value = 25

# Row 21 - 4516351202_UTL
$ws.Range("C21").Value = 0
$ws.Range("G21").Value = 36
$ws.Range("I21").Value = 0

# Row 24 - 47225672
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("I24").Value = -1
